$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-22, column C (ShipmentTrackNum) get new tracking numbers.
# Rows 5,6,7,13,14,15,16,17 also mirror the same value into column D (PackageTrackNum).

$updates = @(
    @{ Row = 2;  Value = "320018655634"; MirrorD = $false },
    @{ Row = 3;  Value = "320018655645"; MirrorD = $false },
    @{ Row = 4;  Value = "320018655678"; MirrorD = $false },
    @{ Row = 5;  Value = "320018655690"; MirrorD = $true  },
    @{ Row = 6;  Value = "320018655760"; MirrorD = $true  },
    @{ Row = 7;  Value = "320018655781"; MirrorD = $true  },
    @{ Row = 8;  Value = "320018655818"; MirrorD = $false },
    @{ Row = 9;  Value = "320018655830"; MirrorD = $false },
    @{ Row = 10; Value = "320018655862"; MirrorD = $false },
    @{ Row = 11; Value = "320018655884"; MirrorD = $false },
    @{ Row = 12; Value = "320018655921"; MirrorD = $false },
    @{ Row = 13; Value = "320018655943"; MirrorD = $true  },
    @{ Row = 14; Value = "320018655976"; MirrorD = $true  },
    @{ Row = 15; Value = "320018655998"; MirrorD = $true  },
    @{ Row = 16; Value = "320018645839"; MirrorD = $true  },
    @{ Row = 17; Value = "320018645850"; MirrorD = $true  },
    @{ Row = 18; Value = "320018645894"; MirrorD = $false },
    @{ Row = 19; Value = "320018645910"; MirrorD = $false },
    @{ Row = 20; Value = "320018645942"; MirrorD = $false },
    @{ Row = 21; Value = "320018645964"; MirrorD = $false },
    @{ Row = 22; Value = "320018645997"; MirrorD = $false }
)

foreach ($u in $updates) {
    $cellC = $ws.Cells.Item($u.Row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $u.Value
    if ($u.MirrorD) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.Value
    }
}
